$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.312497333333334
$ws.Range("H2").Value = 21.937492
$ws.Range("I2").Value = 0.05970572560549242
$ws.Range("J2").Value = 0.05970572560549242
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.722831
$ws.Range("N2").Value = 11.168493
$ws.Range("O2").Value = 0.2042994277655142
$ws.Range("P2").Value = 0.2042994277655142
$ws.Range("Q2").Value = 27.22319175995067
$ws.Range("R2").Value = 245.008725839556
$ws.Range("S2").Value = 0.01219784557552691
$ws.Range("T2").Value = 0.01219784557552691

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.312497333333334
$ws.Range("H3").Value = 21.937492
$ws.Range("I3").Value = 0.05970572560549242
$ws.Range("J3").Value = 0.05970572560549242
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.746044
$ws.Range("N3").Value = 5.238131999999999
$ws.Range("O3").Value = 0.09581842153280916
$ws.Range("P3").Value = 0.09581842153280916
$ws.Range("Q3").Value = 12.76794209388266
$ws.Range("R3").Value = 114.911478844944
$ws.Range("S3").Value = 0.00572090838398931
$ws.Range("T3").Value = 0.00572090838398931

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.312497333333334
$ws.Range("H4").Value = 21.937492
$ws.Range("I4").Value = 0.05970572560549242
$ws.Range("J4").Value = 0.05970572560549242
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 11.46242466666667
$ws.Range("N4").Value = 34.387274
$ws.Range("O4").Value = 0.6290285001401661
$ws.Range("P4").Value = 0.6290285001401661
$ws.Range("Q4").Value = 83.81894980853421
$ws.Range("R4").Value = 754.370548276808
$ws.Range("S4").Value = 0.03755660302740321
$ws.Range("T4").Value = 0.03755660302740321

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.312497333333334
$ws.Range("H5").Value = 21.937492
$ws.Range("I5").Value = 0.05970572560549242
$ws.Range("J5").Value = 0.05970572560549242
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.291125333333333
$ws.Range("N5").Value = 3.873376
$ws.Range("O5").Value = 0.07085365056151052
$ws.Range("P5").Value = 0.07085365056151052
$ws.Range("Q5").Value = 9.441350556999112
$ws.Range("R5").Value = 84.972155012992
$ws.Range("S5").Value = 0.004230368618572991
$ws.Range("T5").Value = 0.004230368618572991

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.68221266666667
$ws.Range("H6").Value = 89.046638
$ws.Range("I6").Value = 0.242351957758873
$ws.Range("J6").Value = 0.242351957758873
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.722831
$ws.Range("N6").Value = 11.168493
$ws.Range("O6").Value = 0.2042994277655142
$ws.Range("P6").Value = 0.2042994277655142
$ws.Range("Q6").Value = 110.5018614640593
$ws.Range("R6").Value = 994.5167531765341
$ws.Range("S6").Value = 0.04951236628798983
$ws.Range("T6").Value = 0.04951236628798982

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.68221266666667
$ws.Range("H7").Value = 89.046638
$ws.Range("I7").Value = 0.242351957758873
$ws.Range("J7").Value = 0.242351957758873
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.746044
$ws.Range("N7").Value = 5.238131999999999
$ws.Range("O7").Value = 0.09581842153280916
$ws.Range("P7").Value = 0.09581842153280916
$ws.Range("Q7").Value = 51.82644933335732
$ws.Range("R7").Value = 466.4380440002159
$ws.Range("S7").Value = 0.02322178204784125
$ws.Range("T7").Value = 0.02322178204784125

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 29.68221266666667
$ws.Range("H8").Value = 89.046638
$ws.Range("I8").Value = 0.242351957758873
$ws.Range("J8").Value = 0.242351957758873
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 11.46242466666667
$ws.Range("N8").Value = 34.387274
$ws.Range("O8").Value = 0.6290285001401661
$ws.Range("P8").Value = 0.6290285001401661
$ws.Range("Q8").Value = 340.2301266316458
$ws.Range("R8").Value = 3062.071139684812
$ws.Range("S8").Value = 0.1524462884950968
$ws.Range("T8").Value = 0.1524462884950968

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 29.68221266666667
$ws.Range("H9").Value = 89.046638
$ws.Range("I9").Value = 0.242351957758873
$ws.Range("J9").Value = 0.242351957758873
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.291125333333333
$ws.Range("N9").Value = 3.873376
$ws.Range("O9").Value = 0.07085365056151052
$ws.Range("P9").Value = 0.07085365056151052
$ws.Range("Q9").Value = 38.32345672332089
$ws.Range("R9").Value = 344.911110509888
$ws.Range("S9").Value = 0.01717152092794515
$ws.Range("T9").Value = 0.01717152092794514

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.06470466666667
$ws.Range("H10").Value = 33.194114
$ws.Range("I10").Value = 0.09034208022509747
$ws.Range("J10").Value = 0.09034208022509749
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.722831
$ws.Range("N10").Value = 11.168493
$ws.Range("O10").Value = 0.2042994277655142
$ws.Range("P10").Value = 0.2042994277655142
$ws.Range("Q10").Value = 41.19202553891133
$ws.Range("R10").Value = 370.728229850202
$ws.Range("S10").Value = 0.01845683529313359
$ws.Range("T10").Value = 0.0184568352931336

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 11.06470466666667
$ws.Range("H11").Value = 33.194114
$ws.Range("I11").Value = 0.09034208022509747
$ws.Range("J11").Value = 0.09034208022509749
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.746044
$ws.Range("N11").Value = 5.238131999999999
$ws.Range("O11").Value = 0.09581842153280916
$ws.Range("P11").Value = 0.09581842153280916
$ws.Range("Q11").Value = 19.31946119500533
$ws.Range("R11").Value = 173.875150755048
$ws.Range("S11").Value = 0.008656435525159252
$ws.Range("T11").Value = 0.008656435525159253

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 11.06470466666667
$ws.Range("H12").Value = 33.194114
$ws.Range("I12").Value = 0.09034208022509747
$ws.Range("J12").Value = 0.09034208022509749
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 11.46242466666667
$ws.Range("N12").Value = 34.387274
$ws.Range("O12").Value = 0.6290285001401661
$ws.Range("P12").Value = 0.6290285001401661
$ws.Range("Q12").Value = 126.8283437005818
$ws.Range("R12").Value = 1141.455093305236
$ws.Range("S12").Value = 0.05682774322353562
$ws.Range("T12").Value = 0.05682774322353563

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 11.06470466666667
$ws.Range("H13").Value = 33.194114
$ws.Range("I13").Value = 0.09034208022509747
$ws.Range("J13").Value = 0.09034208022509749
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.291125333333333
$ws.Range("N13").Value = 3.873376
$ws.Range("O13").Value = 0.07085365056151052
$ws.Range("P13").Value = 0.07085365056151052
$ws.Range("Q13").Value = 14.28592050098489
$ws.Range("R13").Value = 128.573284508864
$ws.Range("S13").Value = 0.006401066183269006
$ws.Range("T13").Value = 0.006401066183269007

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 74.41623166666666
$ws.Range("H14").Value = 223.248695
$ws.Range("I14").Value = 0.6076002364105371
$ws.Range("J14").Value = 0.6076002364105371
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.722831
$ws.Range("N14").Value = 11.168493
$ws.Range("O14").Value = 0.2042994277655142
$ws.Range("P14").Value = 0.2042994277655142
$ws.Range("Q14").Value = 277.0390541518483
$ws.Range("R14").Value = 2493.351487366635
$ws.Range("S14").Value = 0.1241323806088639
$ws.Range("T14").Value = 0.1241323806088639

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 74.41623166666666
$ws.Range("H15").Value = 223.248695
$ws.Range("I15").Value = 0.6076002364105371
$ws.Range("J15").Value = 0.6076002364105371
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.746044
$ws.Range("N15").Value = 5.238131999999999
$ws.Range("O15").Value = 0.09581842153280916
$ws.Range("P15").Value = 0.09581842153280916
$ws.Range("Q15").Value = 129.9340148041933
$ws.Range("R15").Value = 1169.40613323774
$ws.Range("S15").Value = 0.05821929557581934
$ws.Range("T15").Value = 0.05821929557581934

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 74.41623166666666
$ws.Range("H16").Value = 223.248695
$ws.Range("I16").Value = 0.6076002364105371
$ws.Range("J16").Value = 0.6076002364105371
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 11.46242466666667
$ws.Range("N16").Value = 34.387274
$ws.Range("O16").Value = 0.6290285001401661
$ws.Range("P16").Value = 0.6290285001401661
$ws.Range("Q16").Value = 852.9904494563809
$ws.Range("R16").Value = 7676.914045107429
$ws.Range("S16").Value = 0.3821978653941305
$ws.Range("T16").Value = 0.3821978653941305

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 74.41623166666666
$ws.Range("H17").Value = 223.248695
$ws.Range("I17").Value = 0.6076002364105371
$ws.Range("J17").Value = 0.6076002364105371
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.291125333333333
$ws.Range("N17").Value = 3.873376
$ws.Range("O17").Value = 0.07085365056151052
$ws.Range("P17").Value = 0.07085365056151052
$ws.Range("Q17").Value = 96.08068191603554
$ws.Range("R17").Value = 864.7261372443199
$ws.Range("S17").Value = 0.04305069483172338
$ws.Range("T17").Value = 0.04305069483172338
